$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "65.711.84"

# Row 3
Set-TextValue "D3" "2.668.72"
Set-TextValue "E3" "  -1.01%  "

# Row 4
Set-TextValue "E4" "  +0.00%  "

# Row 5
Set-TextValue "D5" "597.78"
Set-TextValue "E5" "  -2.31%  "

# Row 6
Set-TextValue "D6" "157.00"
Set-TextValue "E6" "  -0.93%  "

# Row 7
Set-TextValue "E7" "  +0.00%  "

# Row 8
Set-TextValue "D8" "0.615"
Set-TextValue "E8" "  +4.18%  "

# Row 9
Set-TextValue "E9" "  +2.49%  "

# Row 10
Set-TextValue "D10" "0.399"
Set-TextValue "E10" "  -1.00%  "

# Row 11
Set-TextValue "E11" "  -3.62%  "

# Row 12
Set-TextValue "E12" "  -0.15%  "

# Row 13
Set-TextValue "D13" "29.03"
Set-TextValue "E13" "  -3.66%  "

# Row 14
Set-TextValue "E14" "  -4.95%  "

# Row 15
Set-TextValue "D15" "3.146.36"
Set-TextValue "E15" "  -1.06%  "

# Row 16
Set-TextValue "D16" "65.581.37"
Set-TextValue "E16" "  -0.48%  "

# Row 17
Set-TextValue "D17" "2.671.55"
Set-TextValue "E17" "  -0.78%  "

# Row 18
Set-TextValue "E18" "  -0.72%  "

# Row 19
Set-TextValue "D19" "4.78"
Set-TextValue "E19" "  -2.53%  "

# Row 20
Set-TextValue "D20" "351.69"
Set-TextValue "E20" "  -2.10%  "

# Row 21
Set-TextValue "D21" "7.48"
Set-TextValue "E21" "  -3.89%  "

# Row 22
Set-TextValue "E22" "  -0.06%  "

# Row 23
Set-TextValue "D23" "69.15"
Set-TextValue "E23" "  -3.15%  "

# Row 24
Set-TextValue "E24" "  -1.23%  "

# Row 25
Set-TextValue "E25" "  -3.15%  "

# Row 26
Set-TextValue "D26" "1.68"
Set-TextValue "E26" "  +2.86%  "

# Row 27
Set-TextValue "E27" "  -3.75%  "

# Row 28
Set-TextValue "E28" "  -3.85%  "

# Row 29
Set-TextValue "D29" "7.99"
Set-TextValue "E29" "  -3.63%  "

# Row 30
Set-TextValue "E30" "  -0.03%  "

# Row 31
Set-TextValue "D31" "534.66"
Set-TextValue "E31" "  +0.10%  "

# Row 32
Set-TextValue "E32" "  -3.92%  "

# Row 33
Set-TextValue "E33" "  -1.36%  "

# Row 34
Set-TextValue "E34" "  -3.85%  "

# Row 35
Set-TextValue "D35" "5.47"
Set-TextValue "E35" "  -0.27%  "

# Row 36
Set-TextValue "E36" "  -2.72%  "

# Row 37
Set-TextValue "D37" "20.52"
Set-TextValue "E37" "  -1.19%  "

# Row 38
Set-TextValue "E38" "  -0.02%  "

# Row 39
Set-TextValue "D39" "156.73"
Set-TextValue "E39" "  -3.29%  "

# Row 40
Set-TextValue "E40" "  -3.20%  "

# Row 41
Set-TextValue "D41" "0.999"
Set-TextValue "E41" "  +0.01%  "

# Row 42
Set-TextValue "D42" "162.48"
Set-TextValue "E42" "  -3.36%  "

# Row 43
Set-TextValue "E43" "  -1.75%  "

# Row 44
Set-TextValue "E44" "  +0.67%  "

# Row 45
Set-TextValue "D45" "0.0610"
Set-TextValue "E45" "  -3.97%  "

# Row 46
Set-TextValue "D46" "22.59"
Set-TextValue "E46" "  -5.22%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D47" "0.638"
Set-TextValue "E47" "  -3.06%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D48" "0.0257"
Set-TextValue "E48" "  -4.22%  "

# Row 49
Set-TextValue "E49" "  +9.02%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D50" "0.0987"
Set-TextValue "E50" "  -1.09%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "19.93"
Set-TextValue "E51" "  -4.85%  "
